# Update crypto price/volume table with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.561.72"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.921.96"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4798"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2887"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06719"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "104.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07743"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.921.55"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6863"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "266.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.88%  "
$ws.Range("D17").Value = "30.594.24"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007541"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.460"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.369"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.688"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.099"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1027"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.387"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.674"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.521"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.279"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04776"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7439"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.639"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.378"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8658"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4317"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.611"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.24%  "
$ws.Range("D48").Value = "1.000.65"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.008"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.66%  "
